$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 9).Value = 0.07965462187692204
$ws.Cells.Item(2, 10).Value = 0.07965462187692204
$ws.Cells.Item(2, 16).Value = 0.9810128591839974
$ws.Cells.Item(2, 19).Value = 0.0781422083546995
$ws.Cells.Item(2, 20).Value = 0.07814220835469948

# Row 3
$ws.Cells.Item(3, 9).Value = 0.07965462187692204
$ws.Cells.Item(3, 10).Value = 0.07965462187692204
$ws.Cells.Item(3, 19).Value = 0.001512413522222555
$ws.Cells.Item(3, 20).Value = 0.001512413522222555

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.551351
$ws.Cells.Item(4, 8).Value = 1.654053
$ws.Cells.Item(4, 9).Value = 0.5089207502863742
$ws.Cells.Item(4, 10).Value = 0.5089207502863742
$ws.Cells.Item(4, 16).Value = 0.9810128591839974
$ws.Cells.Item(4, 17).Value = 0.005678363949000001
$ws.Cells.Item(4, 18).Value = 0.051105275541
$ws.Cells.Item(4, 19).Value = 0.4992578003365011
$ws.Cells.Item(4, 20).Value = 0.4992578003365011

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.551351
$ws.Cells.Item(5, 8).Value = 1.654053
$ws.Cells.Item(5, 9).Value = 0.5089207502863742
$ws.Cells.Item(5, 10).Value = 0.5089207502863742
$ws.Cells.Item(5, 17).Value = 0.0001099026326666667
$ws.Cells.Item(5, 18).Value = 0.000989123694
$ws.Cells.Item(5, 19).Value = 0.009662949949873051
$ws.Cells.Item(5, 20).Value = 0.009662949949873051

# Row 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.018049
$ws.Cells.Item(6, 8).Value = 0.054147
$ws.Cells.Item(6, 9).Value = 0.01666000537211099
$ws.Cells.Item(6, 10).Value = 0.01666000537211099
$ws.Cells.Item(6, 16).Value = 0.9810128591839974
$ws.Cells.Item(6, 17).Value = 0.000185886651
$ws.Cells.Item(6, 18).Value = 0.001672979859
$ws.Cells.Item(6, 19).Value = 0.01634367950411536
$ws.Cells.Item(6, 20).Value = 0.01634367950411536

# Row 7
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.018049
$ws.Cells.Item(7, 8).Value = 0.054147
$ws.Cells.Item(7, 9).Value = 0.01666000537211099
$ws.Cells.Item(7, 10).Value = 0.01666000537211099
$ws.Cells.Item(7, 17).Value = 0.000003597767333333333
$ws.Cells.Item(7, 18).Value = 0.000032379906
$ws.Cells.Item(7, 19).Value = 0.0003163258679956301
$ws.Cells.Item(7, 20).Value = 0.0003163258679956302

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.4276773333333333
$ws.Cells.Item(8, 8).Value = 1.283032
$ws.Cells.Item(8, 9).Value = 0.3947646224645928
$ws.Cells.Item(8, 10).Value = 0.3947646224645929
$ws.Cells.Item(8, 16).Value = 0.9810128591839974
$ws.Cells.Item(8, 17).Value = 0.004404648856
$ws.Cells.Item(8, 18).Value = 0.039641839704
$ws.Cells.Item(8, 19).Value = 0.3872691709886815
$ws.Cells.Item(8, 20).Value = 0.3872691709886815

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.4276773333333333
$ws.Cells.Item(9, 8).Value = 1.283032
$ws.Cells.Item(9, 9).Value = 0.3947646224645928
$ws.Cells.Item(9, 10).Value = 0.3947646224645929
$ws.Cells.Item(9, 17).Value = 0.00008525034844444444
$ws.Cells.Item(9, 18).Value = 0.000767253136
$ws.Cells.Item(9, 19).Value = 0.007495451475911304
$ws.Cells.Item(9, 20).Value = 0.007495451475911304
